# Update gh-pages to output generated at 456a3b4
# Applies numeric "want-to-go counter" bumps across all four sheets, a
# sold-out status flip, a local-life entry that expired (row removed /
# everything below it shifts up), and a local-life entry that got replaced
# by a new one (content swap) on "全部类型".

$wb = $excel.ActiveWorkbook

function Set-Num($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-Txt($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Text that LOOKS like a bare date ("2024-07-26") gets auto-parsed into a
# real date serial by the Value setter. Force it to stay literal text by
# quote-prefixing, then strip the resulting style back to Normal so no
# stray number-format sticks to the cell.
function Set-DateTxt($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) -- plain "想去人数" (want-to-go count) bumps
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$s1F = @{
    3 = 255; 4 = 563; 5 = 2479; 7 = 149; 10 = 4979; 11 = 6238; 12 = 892;
    13 = 85; 15 = 1345; 16 = 556; 17 = 6803; 18 = 379; 19 = 30; 21 = 4580;
    22 = 377; 25 = 2230; 27 = 416; 28 = 1131; 29 = 180; 31 = 65; 32 = 132;
    35 = 1958; 36 = 199; 37 = 482; 39 = 1326; 40 = 583; 42 = 4; 43 = 86;
    44 = 1063; 45 = 1346; 48 = 218; 49 = 60
}
foreach ($r in $s1F.Keys) { Set-Num $ws1 $r 6 $s1F[$r] }

# Row 22 also flipped from a numeric lowest-price to "已售罄" (sold out).
Set-Txt $ws1 22 7 "已售罄"

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performance) -- plain want-to-go count bumps
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$s2F = @{
    4 = 424; 15 = 24; 16 = 159; 19 = 233; 20 = 125; 23 = 138; 27 = 253; 28 = 20
}
foreach ($r in $s2F.Keys) { Set-Num $ws2 $r 6 $s2F[$r] }

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life)
#   - F6/F7 plain bumps
#   - Row 8 ("剧场版BLUE LOCK" cafe, now expired) is removed; rows 9-14
#     shift up into rows 8-13 (with their own want-to-go counts bumped),
#     and the now-unused row 14 disappears (dimension A1:I14 -> A1:I13).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-Num $ws3 6 6 1627
Set-Num $ws3 7 6 525

# Capture old rows 9..14 (B..I) before they get overwritten.
$cols = 2..9
$oldRows = @{}
foreach ($r in 9..14) {
    $vals = @{}
    foreach ($c in $cols) { $vals[$c] = $ws3.Cells.Item($r, $c).Value() }
    $oldRows[$r] = $vals
}

# New F (want-to-go count) values for the shifted rows, keyed by the OLD
# row number whose data is moving (old r -> new r-1).
$shiftedF = @{ 9 = 1206; 10 = 1196; 11 = 1700; 12 = 2048; 13 = 514; 14 = 419 }

foreach ($oldR in 9..14) {
    $newR = $oldR - 1
    $vals = $oldRows[$oldR]
    Set-DateTxt $ws3 $newR 2 $vals[2]          # B: date
    Set-Txt     $ws3 $newR 3 $vals[3]          # C: name
    Set-Txt     $ws3 $newR 4 $vals[4]          # D: location
    Set-Txt     $ws3 $newR 5 $vals[5]          # E: time range
    Set-Num     $ws3 $newR 6 $shiftedF[$oldR]  # F: want-to-go count (bumped)
    $g = $vals[7]
    if ($g -is [string]) { Set-Txt $ws3 $newR 7 $g } else { Set-Num $ws3 $newR 7 $g }  # G
    Set-Txt $ws3 $newR 8 $vals[8]               # H: link
    Set-Txt $ws3 $newR 9 $vals[9]               # I: cover image
}

# Row 14 is now a duplicate of (new) row 13; drop it entirely so the sheet
# shrinks back to 13 data rows, matching dimension A1:I13.
$ws3.Rows.Item(14).Delete()

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All Types) -- mirrors the other sheets' bumps, plus
# its own two special edits.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$s4F = @{
    2 = 1627; 3 = 563; 4 = 525; 5 = 2479; 6 = 1206; 8 = 2048; 9 = 4979;
    10 = 514; 14 = 85; 17 = 1345; 18 = 556; 19 = 6803; 20 = 379; 21 = 419;
    24 = 4580; 27 = 2230; 29 = 416; 30 = 1131; 31 = 180; 32 = 65; 35 = 132;
    37 = 1958; 38 = 199; 39 = 482; 42 = 1326; 43 = 138; 44 = 583; 46 = 20;
    47 = 1063; 48 = 1346; 49 = 218
}
foreach ($r in $s4F.Keys) { Set-Num $ws4 $r 6 $s4F[$r] }

# Row 13's count wasn't bumped -- it was reset to 0.
Set-Num $ws4 13 6 0

# Row 25: the old "创造力动漫游戏嘉年华签售票-爱拍照的玉老师" listing was
# swapped out entirely for a new "动漫水着嘉年华" listing (date B25 stays
# "2024-08-10", everything else changes).
Set-Txt $ws4 25 3 "上海·动漫水着嘉年华"
Set-Txt $ws4 25 4 "民府路678号 抖音江湾广场"
Set-Txt $ws4 25 5 "2024.08.10 11:00-08.10 17:00"
Set-Num $ws4 25 6 29
Set-Num $ws4 25 7 68
Set-Txt $ws4 25 8 "https://show.bilibili.com/platform/detail.html?id=89929"
Set-Txt $ws4 25 9 "//i2.hdslb.com/bfs/openplatform/202407/UTF8WBkE1722219440039.jpeg"
